# Fixed #476 Moving from Apache POI 4.1.0 to 5.2.3.
#
# The three list items ("Coffee", "Tea", "Milk") each carry an explicit
# run-level Bold / Italic / Strikethrough = Off formatting. Re-assert that
# formatting on each run through the Word object model so the run
# properties are rewritten (matching the newer POI generator's output for
# the same semantic "false" toggle state).

$d = $word.ActiveDocument
$words = @("Coffee", "Tea", "Milk")

foreach ($w in $words) {
    $rng = $d.Content
    $found = $rng.Find.Execute($w, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = $false
        $rng.Font.Italic = $false
        $rng.Font.StrikeThrough = $false
    }
}
